$d = $word.ActiveDocument

$replacements = @(
    @("96×96=9216", "71×32=2272"),
    @("32×71=2272", "97×48=4656"),
    @("94×55=5170", "49×50=2450"),
    @("18×94=1692", "69×23=1587"),
    @("57×49=2793", "42×52=2184"),
    @("96×80=7680", "63×72=4536"),
    @("39×30=1170", "81×15=1215"),
    @("84×99=8316", "42×51=2142"),
    @("86×66=5676", "80×56=4480"),
    @("17×83=1411", "29×50=1450"),
    @("69×86=5934", "37×23=851"),
    @("53×48=2544", "38×44=1672"),
    @("20×78=1560", "20×69=1380"),
    @("48×64=3072", "83×34=2822"),
    @("32×47=1504", "89×94=8366"),
    @("87×90=7830", "74×89=6586"),
    @("36×93=3348", "27×92=2484"),
    @("77×94=7238", "18×18=324"),
    @("15×30=450",  "77×52=4004"),
    @("13×15=195",  "22×48=1056"),
    @("56×64=3584", "99×33=3267"),
    @("94×16=1504", "51×96=4896"),
    @("54×34=1836", "45×62=2790"),
    @("76×58=4408", "58×73=4234"),
    @("40×79=3160", "56×73=4088")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done replacing $($replacements.Count) values"
